# Refresh currentAveragePrice/Leve price/profit figures pulled by the scheduled
# market-data runner. Each block updates one leve row's H:N price/profit columns
# to the latest fetched values (plain values, no formulas are used in this sheet).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 15: Morning Glass of Ether
$ws.Range("H15").Value = 1385.3855
$ws.Range("I15").Value = 1385.3855
$ws.Range("K15").Value = 4156.1565
$ws.Range("M15").Value = -3987.1565

# Row 33: Glazed and Confused
$ws.Range("H33").Value = 117.30769
$ws.Range("I33").Value = 117.30769
$ws.Range("K33").Value = 117.30769
$ws.Range("M33").Value = 111.69231

# Row 43: Growing Is Knowing
$ws.Range("H43").Value = 1328.1428
$ws.Range("I43").Value = 1399.75
$ws.Range("J43").Value = 1299.5
$ws.Range("K43").Value = 1399.75
$ws.Range("L43").Value = 1299.5
$ws.Range("M43").Value = -1330.75
$ws.Range("N43").Value = -1437.5

# Row 138: All-night Crafting
$ws.Range("H138").Value = 1742.55
$ws.Range("I138").Value = 918.4103
$ws.Range("J138").Value = 2269.459
$ws.Range("K138").Value = 2755.2309
$ws.Range("L138").Value = 6808.376999999999
$ws.Range("M138").Value = 2384.7691
$ws.Range("N138").Value = -17088.377

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust
$ws.Range("H32").Value = 4585.2
$ws.Range("I32").Value = 2696.3037
$ws.Range("J32").Value = 13911.625
$ws.Range("K32").Value = 2696.3037
$ws.Range("L32").Value = 13911.625
$ws.Range("M32").Value = -2409.3037
$ws.Range("N32").Value = -14485.625

# Row 97: Ore for Me
$ws.Range("H97").Value = 1098.1428
$ws.Range("I97").Value = 818.4211
$ws.Range("K97").Value = 818.4211
$ws.Range("M97").Value = -322.4211

# Row 132: Don't Bore Me, Ore Me
$ws.Range("H132").Value = 1442.9459
$ws.Range("I132").Value = 1187.8163
$ws.Range("J132").Value = 1943
$ws.Range("K132").Value = 3563.4489
$ws.Range("L132").Value = 5829
$ws.Range("M132").Value = -1033.4489
$ws.Range("N132").Value = -10889

$ws = $wb.Worksheets.Item("BSM")
# Row 94: High Steal
$ws.Range("H94").Value = 1026.4706
$ws.Range("I94").Value = 496.73334
$ws.Range("K94").Value = 496.73334
$ws.Range("M94").Value = -45.73334

# Row 135: Axes to the Maxes
$ws.Range("H135").Value = 59600
$ws.Range("J135").Value = 59600
$ws.Range("L135").Value = 59600
$ws.Range("N135").Value = -69740

$ws = $wb.Worksheets.Item("CRP")
# Row 10: Spears and Sorcery
$ws.Range("H10").Value = 1044.5555
$ws.Range("I10").Value = 800
$ws.Range("J10").Value = 3001
$ws.Range("K10").Value = 800
$ws.Range("L10").Value = 3001
$ws.Range("M10").Value = -661
$ws.Range("N10").Value = -3279

# Row 132: Hull Lotta Damage
$ws.Range("H132").Value = 1495.1333
$ws.Range("I132").Value = 1411
$ws.Range("J132").Value = 1621.3334
$ws.Range("K132").Value = 4233
$ws.Range("L132").Value = 4864.0002
$ws.Range("M132").Value = -1703
$ws.Range("N132").Value = -9924.0002

# Row 140: Spear Pressure
$ws.Range("H140").Value = 58000
$ws.Range("J140").Value = 58000
$ws.Range("L140").Value = 58000
$ws.Range("N140").Value = -68360

$ws = $wb.Worksheets.Item("CUL")
# Row 23: Sweet Smell of Success
$ws.Range("H23").Value = 216.46153
$ws.Range("I23").Value = 105.5
$ws.Range("J23").Value = 394
$ws.Range("K23").Value = 316.5
$ws.Range("L23").Value = 1182
$ws.Range("M23").Value = -81.5
$ws.Range("N23").Value = -1652

# Row 115: Mixology
$ws.Range("H115").Value = 4339.467
$ws.Range("I115").Value = 1519.8
$ws.Range("J115").Value = 5749.3
$ws.Range("K115").Value = 4559.4
$ws.Range("L115").Value = 17247.9
$ws.Range("M115").Value = -3384.4
$ws.Range("N115").Value = -19597.9

# Row 131: The Mountain Steeped
$ws.Range("H131").Value = 18659.059
$ws.Range("I131").Value = 356.33334
$ws.Range("K131").Value = 1069.00002
$ws.Range("M131").Value = 3970.99998

# Row 136: Simple Is Hardest
$ws.Range("H136").Value = 1595.5264
$ws.Range("I136").Value = 1312.6471
$ws.Range("J136").Value = 4000
$ws.Range("K136").Value = 3937.9413
$ws.Range("L136").Value = 12000
$ws.Range("M136").Value = 1162.0587
$ws.Range("N136").Value = -22200

$ws = $wb.Worksheets.Item("GSM")
# Row 70: Sky Is the Limit
$ws.Range("H70").Value = 4435.8887
$ws.Range("I70").Value = 4150.75
$ws.Range("K70").Value = 4150.75
$ws.Range("M70").Value = -3880.75

# Row 73: Hulls of Broken Dreams (L)
$ws.Range("H73").Value = 4435.8887
$ws.Range("I73").Value = 4150.75
$ws.Range("K73").Value = 4150.75
$ws.Range("M73").Value = -3214.75

# Row 107: Whetstones for the Workers
$ws.Range("H107").Value = 163.42857
$ws.Range("I107").Value = 163.42857
$ws.Range("K107").Value = 163.42857
$ws.Range("M107").Value = 1756.57143

# Row 132: On Board for Lar
$ws.Range("H132").Value = 1834946.9
$ws.Range("I132").Value = 2567177.2
$ws.Range("J132").Value = 4370.8335
$ws.Range("K132").Value = 7701531.600000001
$ws.Range("L132").Value = 13112.5005
$ws.Range("M132").Value = -7699001.600000001
$ws.Range("N132").Value = -18172.5005

$ws = $wb.Worksheets.Item("LTW")
# Row 42: Slave to Fashion
$ws.Range("H42").Value = 24750
$ws.Range("J42").Value = 24750
$ws.Range("L42").Value = 24750
$ws.Range("N42").Value = -25876

# Row 46: Supply Side Logic
$ws.Range("H46").Value = 1878
$ws.Range("I46").Value = 1343.625
$ws.Range("K46").Value = 1343.625
$ws.Range("M46").Value = -1155.625

# Row 49: First They Came for the Heretics
$ws.Range("H49").Value = 24750
$ws.Range("J49").Value = 24750
$ws.Range("L49").Value = 24750
$ws.Range("N49").Value = -25044

# Row 55: It's Not a Job, It's a Calling
$ws.Range("H55").Value = 453.4091
$ws.Range("I55").Value = 348.3846
$ws.Range("J55").Value = 605.1111
$ws.Range("K55").Value = 348.3846
$ws.Range("L55").Value = 605.1111
$ws.Range("M55").Value = -175.3846
$ws.Range("N55").Value = -951.1111

# Row 104: Brace Yourselves
$ws.Range("H104").Value = 200000
$ws.Range("J104").Value = 200000
$ws.Range("L104").Value = 200000
$ws.Range("N104").Value = -206988

# Row 132: Tenets of Tanning
$ws.Range("H132").Value = 2258.186
$ws.Range("I132").Value = 1970.3334
$ws.Range("K132").Value = 5911.0002
$ws.Range("M132").Value = -3381.0002

$ws = $wb.Worksheets.Item("WVR")
# Row 96: Skills on Display
$ws.Range("H96").Value = 5472.2104
$ws.Range("J96").Value = 9774
$ws.Range("L96").Value = 9774
$ws.Range("N96").Value = -12520

# Row 100: Of Great Import
$ws.Range("H100").Value = 984.9231
$ws.Range("I100").Value = 813.44446
$ws.Range("K100").Value = 1626.88892
$ws.Range("M100").Value = -1085.88892

# Row 122: Heavy Armoire
$ws.Range("H122").Value = 33699.16
$ws.Range("I122").Value = 39826.81
$ws.Range("J122").Value = 1529
$ws.Range("K122").Value = 119480.43
$ws.Range("L122").Value = 4587
$ws.Range("M122").Value = -117030.43
$ws.Range("N122").Value = -9487

# Row 132: Comfy Cabins
$ws.Range("H132").Value = 1464.9062
$ws.Range("I132").Value = 1671.9048
$ws.Range("J132").Value = 1069.7273
$ws.Range("K132").Value = 5015.7144
$ws.Range("L132").Value = 3209.1819
$ws.Range("M132").Value = -2485.7144
$ws.Range("N132").Value = -8269.1819
